# Applies the edits described by the commit diff:
#   - Restricciones_del_follower: new coefficients/bounds for J_0_L0_v,
#     J_0_LP_v and J_Ne_L0_v rows
#   - Punto_modificado, Vector_bf, Vector_BF: refreshed x/y point and
#     derived vectors
#   - Vector_Alpha: updated alpha scalar
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Restricciones_del_follower")
$ws.Range("A2").Value = "-15.454946236559138 - x + 3.2795698924731185y"

$ws.Range("B2").Value = "'17.454946236559138"
$ws.Range("B2").Style = "Normal"

$ws.Range("D2").Value = "'0.93"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "'4.5"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = "'6.1"
$ws.Range("F2").Style = "Normal"

$ws.Range("A3").Value = "-21.030268817204302 - 0.25x + 3.602150537634409y"

$ws.Range("B3").Value = "'19.030268817204302"
$ws.Range("B3").Style = "Normal"

$ws.Range("D3").Value = "'0.41"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "'4.9"
$ws.Range("E3").Style = "Normal"

$ws.Range("F3").Value = "'6.7"
$ws.Range("F3").Style = "Normal"

$ws.Range("A4").Value = "-4.78 + x"

$ws.Range("B4").Value = "'-3.2199999999999998"
$ws.Range("B4").Style = "Normal"

$ws.Range("D4").Value = "'0.7"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = "'8.5"
$ws.Range("E4").Style = "Normal"

$ws.Range("F4").Value = "'2.9"
$ws.Range("F4").Style = "Normal"

$ws.Range("A5").Value = "-28.19043010752688 + x + 3.763440860215053y"

$ws.Range("B5").Value = "'26.000430107526878"
$ws.Range("B5").Style = "Normal"

$ws.Range("D5").Value = "'0.36"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "'7.3"
$ws.Range("E5").Style = "Normal"

$ws.Range("F5").Value = "'7.0"
$ws.Range("F5").Style = "Normal"

$ws.Range("A6").Value = "-17.55774193548387 + 2.741935483870968y"

$ws.Range("B6").Value = "'16.91774193548387"
$ws.Range("B6").Style = "Normal"

$ws.Range("D6").Value = "'0.79"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "'7.4"
$ws.Range("E6").Style = "Normal"

$ws.Range("F6").Value = "'5.1"
$ws.Range("F6").Style = "Normal"

$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2").Value = "'4.78"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "'6.17"
$ws.Range("B2").Style = "Normal"

$ws = $wb.Worksheets.Item("Vector_bf")
$ws.Range("A2").Value = "'-7.047849462365592"
$ws.Range("A2").Style = "Normal"

# "Vector_bf" and "Vector_BF" only differ by case, and Worksheets.Item
# resolves sheet names case-insensitively (picking Vector_bf first for
# either spelling), so address Vector_BF (the 6th sheet) by its 1-based
# index instead of by name.
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = "'-11.075"
$ws.Range("A2").Style = "Normal"
$ws.Range("A3").Value = "'-81.17204301075269"
$ws.Range("A3").Style = "Normal"

$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 1.8599999999999999
